$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C ("Förändrad") holds a date serial number that was bumped by one day
# (46060 -> 46061) for every data row (rows 2 through 546).
$range = $ws.Range("C2:C546")
$range.Value = 46061
